$d = $word.ActiveDocument

$replacements = @(
    @("2024-03-25 Monday", "2024-03-26 Tuesday"),
    @("91×81=7371", "96×82=7872"),
    @("66×79=5214", "32×74=2368"),
    @("67×56=3752", "85×18=1530"),
    @("57×68=3876", "20×71=1420"),
    @("59×49=2891", "93×13=1209"),
    @("30×34=1020", "96×67=6432"),
    @("54×28=1512", "88×41=3608"),
    @("50×71=3550", "71×61=4331"),
    @("14×59=826", "60×98=5880"),
    @("86×77=6622", "63×34=2142"),
    @("59×25=1475", "11×28=308"),
    @("45×24=1080", "86×25=2150"),
    @("59×45=2655", "62×89=5518"),
    @("21×57=1197", "31×77=2387"),
    @("65×88=5720", "48×63=3024"),
    @("57×62=3534", "44×23=1012"),
    @("60×50=3000", "18×91=1638"),
    @("60×52=3120", "36×63=2268"),
    @("13×69=897", "58×46=2668"),
    @("31×49=1519", "15×24=360"),
    @("63×68=4284", "92×56=5152"),
    @("59×29=1711", "15×99=1485"),
    @("19×50=950", "54×41=2214"),
    @("61×45=2745", "90×18=1620"),
    @("82×86=7052", "29×11=319")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
